# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.980.73'
$ws.Range("E2").Value = '  +1.11%  '

$ws.Range("D3").Value = '2.301.92'
$ws.Range("E3").Value = '  -0.02%  '

$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").Value = '''309.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.12%  '

$ws.Range("D6").Value = '''105.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.29%  '

$ws.Range("E7").Value = '  -0.44%  '

$ws.Range("E8").Value = '  +0.26%  '

$ws.Range("D9").Value = '''0.605'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.09%  '

$ws.Range("D10").Value = '''39.76'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.18%  '

$ws.Range("E11").Value = '  +0.21%  '

$ws.Range("D12").Value = '''8.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.90%  '

$ws.Range("E13").Value = '  +0.13%  '

$ws.Range("D14").Value = '''0.987'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.42%  '

$ws.Range("D15").Value = '''15.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("D16").Value = '2.649.82'
$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("D17").Value = '2.291.70'
$ws.Range("E17").Value = '  -0.17%  '

$ws.Range("D18").Value = '42.830.14'
$ws.Range("E18").Value = '  +0.51%  '

$ws.Range("E19").Value = '  -4.11%  '

$ws.Range("D20").Value = '''13.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.01%  '

$ws.Range("E21").Value = '  -1.32%  '

$ws.Range("D22").Value = '''73.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.84%  '

$ws.Range("D23").Value = '''3.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.78%  '

$ws.Range("D24").Value = '''268.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.22%  '

$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("D27").Value = '''7.58'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.88%  '

$ws.Range("D28").Value = '''10.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.19%  '

$ws.Range("D29").Value = '''2.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.02%  '

$ws.Range("D30").Value = '''37.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.53%  '

$ws.Range("D31").Value = '''22.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.00%  '

$ws.Range("D32").Value = '''165.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.13%  '

$ws.Range("E33").Value = '  -2.08%  '

$ws.Range("E34").Value = '  +6.44%  '

$ws.Range("E35").Value = '  -0.96%  '

$ws.Range("E36").Value = '  -1.09%  '

$ws.Range("D37").Value = '''4.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.31%  '

$ws.Range("E38").Value = '  +0.70%  '

$ws.Range("D39").Value = '''2.79'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.48%  '

$ws.Range("D40").Value = '''3.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.20%  '

$ws.Range("D41").Value = '''107.82'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.21%  '

$ws.Range("E42").Value = '  -3.25%  '

$ws.Range("D43").Value = '''71.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.63%  '

$ws.Range("E44").Value = '  +0.71%  '

$ws.Range("D45").Value = '''1.01'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.31%  '

$ws.Range("D46").Value = '''12.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.94%  '

$ws.Range("D47").Value = '1.703.33'
$ws.Range("E47").Value = '  +3.44%  '

$ws.Range("D48").Value = '''111.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.47%  '

$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").Value = '''75.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.31%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '''8.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.28%  '

$ws.Range("D51").Value = '''5.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.17%  '
